$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Core version (Build) from 7842 to 7853
$ws.Range("C2").Value = 7853

# Add new changelog entries
$ws.Range("A6").Value = "Fixed mipmaps not working"
$ws.Range("A7").Value = "Fixed cloud shaders"

# Move selection to A8, matching the saved workbook state
$ws.Range("A8").Select()
